$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.011.74"
$ws.Range("E2").Value = "  -0.71%  "

# Row 3
$ws.Range("D3").Value = "1.744.84"
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.88"
$ws.Range("E5").Value = "  +2.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.22%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5080"
$ws.Range("E7").Value = "  -6.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2761"
$ws.Range("E8").Value = "  -2.79%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06196"
$ws.Range("E9").Value = "  -0.19%  "

# Row 10
$ws.Range("B10").Value = "TRON"
$ws.Range("C10").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07263"
$ws.Range("E10").Value = "  +0.89%  "

# Row 11
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.741.52"
$ws.Range("E11").Value = "  -1.01%  "

# Row 12
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6535"
$ws.Range("E12").Value = "  -0.82%  "

# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.19"
$ws.Range("E13").Value = "  -2.59%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.693"
$ws.Range("E14").Value = "  +0.66%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.68"
$ws.Range("E15").Value = "  -1.35%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  +0.00%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9995"
$ws.Range("E17").Value = "  -0.19%  "

# Row 18
$ws.Range("D18").Value = "26.017.09"
$ws.Range("E18").Value = "  -0.66%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.90"
$ws.Range("E19").Value = "  -0.66%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006908"
$ws.Range("E20").Value = "  +0.83%  "

# Row 21
$ws.Range("D21").Value = "1.967.24"
$ws.Range("E21").Value = "  -0.84%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.478"
$ws.Range("E22").Value = "  +1.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.760"
$ws.Range("E23").Value = "  -0.66%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.373"
$ws.Range("E24").Value = "  +0.97%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.11"
$ws.Range("E25").Value = "  -3.24%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.509"
$ws.Range("E26").Value = "  -0.46%  "

# Row 27
$ws.Range("E27").Value = "  -0.49%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.782"
$ws.Range("E28").Value = "  -1.89%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.76"
$ws.Range("E29").Value = "  -0.28%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.868"
$ws.Range("E30").Value = "  +1.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08202"
$ws.Range("E31").Value = "  -4.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.650"
$ws.Range("E32").Value = "  -1.45%  "

# Row 33
$ws.Range("E33").Value = "  +0.48%  "

# Row 34
$ws.Range("E34").Value = "  -0.56%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9982"
$ws.Range("E35").Value = "  -1.34%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6127"
$ws.Range("E36").Value = "  -3.15%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.793"
$ws.Range("E37").Value = "  +2.84%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01619"
$ws.Range("E38").Value = "  -0.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.941"
$ws.Range("E39").Value = "  -1.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9998"
$ws.Range("E40").Value = "  -0.25%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.76"
$ws.Range("E41").Value = "  +0.42%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3928"
$ws.Range("E42").Value = "  -0.80%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7697"
$ws.Range("E43").Value = "  +2.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.008"
$ws.Range("E44").Value = "  -1.14%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1161"
$ws.Range("E45").Value = "  +0.32%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.359"
$ws.Range("E46").Value = "  -0.61%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05327"
$ws.Range("E47").Value = "  -0.50%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.57"
$ws.Range("E48").Value = "  +0.60%  "

# Row 49
$ws.Range("E49").Value = "  -1.54%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.610"
$ws.Range("E50").Value = "  -0.37%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3459"
$ws.Range("E51").Value = "  -1.91%  "
